$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing value in row 20 (D20: 0.677 -> 0.678) ---
$ws.Range("D20").Value = 0.678

# --- New row 21: Denoising-AE / noise=20% discussion ---
$ws.Range("A21").Value = "Denoising-AE"
$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = 0.669
$ws.Range("E21").Value = "noise=20%"
$ws.Range("F21").Value = "possible discussion of BN vs no BN"

# --- New row 22: Tied-AE (1/2 params) ---
$ws.Range("A22").Value = "Tied-AE (1/2 params)"
$ws.Range("B22").Value = "-"
$ws.Range("C22").Value = "-"
$ws.Range("D22").Value = 0.611

# --- New row 23: Tied-AE (~#params AE) ---
$ws.Range("A23").Value = "Tied-AE (~#params AE)"
$ws.Range("B23").Value = "-"
$ws.Range("C23").Value = "-"
$ws.Range("D23").Value = 0.659

# --- New row 24: Denoising-Tied-AE ---
$ws.Range("A24").Value = "Denoising-Tied-AE"
$ws.Range("B24").Value = "-"
$ws.Range("C24").Value = "-"
$ws.Range("D24").Value = 0.647

# Re-apply the number format (style index 3, "0.000") to the cells that
# use it, matching D8/D18/D20's styling.
$ws.Range("D22").NumberFormat = $ws.Range("D20").NumberFormat
$ws.Range("D24").NumberFormat = $ws.Range("D20").NumberFormat

# --- Column A widened to fit the new, longer labels ---
$ws.Columns.Item(1).ColumnWidth = 18.75

# --- Update the view: scroll down a bit and move the active selection ---
$ws.Range("E19").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
